# "Generate Report for Archive"
# - Status of "Ready for handoff" rows moves on to "In Translation"
# - The (now narrower) zh-cn / de-de status columns are re-sized to match
#
# Note on the ColumnWidth value: the stored width shrinks from
# 17.2159881591797 to 13.4101845877511 character-width units. ColumnWidth
# is set in characters and gets snapped to the engine's column-width grid,
# so 12.5 is the input that lands on the closest representable width.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: status shown per-locale in columns E (zh-cn) and F (de-de) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status column is C ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status column is C ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
